$wb = $excel.ActiveWorkbook

# Update "想去人数" (column F) values per the diff, across all four sheets.

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 1274
$ws.Cells.Item(5, 6).Value = 1428
$ws.Cells.Item(7, 6).Value = 48
$ws.Cells.Item(8, 6).Value = 684
$ws.Cells.Item(9, 6).Value = 159
$ws.Cells.Item(10, 6).Value = 178
$ws.Cells.Item(13, 6).Value = 17
$ws.Cells.Item(14, 6).Value = 536
$ws.Cells.Item(19, 6).Value = 758
$ws.Cells.Item(20, 6).Value = 2566
$ws.Cells.Item(22, 6).Value = 20
$ws.Cells.Item(25, 6).Value = 184
$ws.Cells.Item(27, 6).Value = 124
$ws.Cells.Item(28, 6).Value = 572
$ws.Cells.Item(31, 6).Value = 131
$ws.Cells.Item(35, 6).Value = 247

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 644
$ws.Cells.Item(9, 6).Value = 290
$ws.Cells.Item(15, 6).Value = 9
$ws.Cells.Item(16, 6).Value = 948
$ws.Cells.Item(23, 6).Value = 255
$ws.Cells.Item(24, 6).Value = 233

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(5, 6).Value = 2331
$ws.Cells.Item(6, 6).Value = 940
$ws.Cells.Item(9, 6).Value = 1180
$ws.Cells.Item(10, 6).Value = 298

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 2331
$ws.Cells.Item(8, 6).Value = 940
$ws.Cells.Item(9, 6).Value = 1180
$ws.Cells.Item(10, 6).Value = 298
$ws.Cells.Item(12, 6).Value = 1274
$ws.Cells.Item(13, 6).Value = 1428
$ws.Cells.Item(15, 6).Value = 48
$ws.Cells.Item(16, 6).Value = 684
$ws.Cells.Item(17, 6).Value = 159
$ws.Cells.Item(19, 6).Value = 178
$ws.Cells.Item(21, 6).Value = 17
$ws.Cells.Item(22, 6).Value = 536
$ws.Cells.Item(26, 6).Value = 758
$ws.Cells.Item(27, 6).Value = 2566
$ws.Cells.Item(29, 6).Value = 20
$ws.Cells.Item(31, 6).Value = 290
$ws.Cells.Item(32, 6).Value = 184
$ws.Cells.Item(33, 6).Value = 124
$ws.Cells.Item(34, 6).Value = 572
$ws.Cells.Item(39, 6).Value = 131
$ws.Cells.Item(40, 6).Value = 9
$ws.Cells.Item(44, 6).Value = 255
$ws.Cells.Item(45, 6).Value = 233
$ws.Cells.Item(49, 6).Value = 247
